# Fixed update to excel issue
# - Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity" sheet
# - Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend" sheet
# - Add a new "PO Forecast" sheet (sheetId 3) with ds / PO_Forecast / yhat_lower / yhat_upper columns

$wb = $excel.ActiveWorkbook

# --- 1. Update existing headers -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row ---------------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the existing header style (bold, bordered, centered) from "Weekly Quantity"
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# --- 4. Data rows ------------------------------------------------------------------
$data = @(
    @(45109.99999999999, 14, -188.858714260272, 213.8227366162305),
    @(45123.99999999999, 44, -158.1055062662614, 229.1782457905653),
    @(45130.99999999999, 59, -146.2485002260766, 243.6244585396699),
    @(45137.99999999999, 74, -108.2546156293941, 271.2687388692803),
    @(45144.99999999999, 88, -100.9559277356119, 285.0360130068379),
    @(45151.99999999999, 103, -93.07577653198234, 308.038007804622),
    @(45158.99999999999, 118, -81.73474403880397, 305.0119520502283),
    @(45172.99999999999, 148, -50.5552212464661, 348.3662611191471),
    @(45179.99999999999, 163, -21.28925371045343, 365.0164580439715),
    @(45186.99999999999, 178, -20.9197227339335, 375.2216858491018),
    @(45193.99999999999, 192, -4.269052105489651, 382.4200197913253),
    @(45200.99999999999, 207, -5.381024233271095, 407.0069994110447),
    @(45249.99999999999, 311, 101.7321139755431, 499.1458015266226),
    @(45256.99999999999, 326, 132.2965326440833, 517.2241847806027),
    @(45270.99999999999, 356, 156.2550372241217, 550.1360639969877),
    @(45277.99999999999, 371, 151.9049238179728, 564.7317049688149),
    @(45298.99999999999, 415, 213.1826002402472, 622.0703144495137),
    @(45305.99999999999, 430, 228.049112920623, 624.9822311758999),
    @(45312.99999999999, 445, 255.3693283485948, 636.900398491346),
    @(45319.99999999999, 460, 248.8865893315744, 664.7457304503345),
    @(45326.99999999999, 475, 280.297783913797, 691.8313905845878),
    @(45333.99999999999, 489, 293.6652108723027, 692.875916868991),
    @(45340.99999999999, 504, 313.9369297346863, 712.1068818704352),
    @(45347.99999999999, 519, 330.9573435869801, 715.4670458489239),
    @(45354.99999999999, 534, 325.682828113774, 736.2643963286756)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    $wsForecast.Cells.Item($rowNum, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowNum, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowNum, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowNum, 4).Value = $row[3]
}

# Reuse the existing date-format style for column A (the "ds" date column)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A" + ($data.Length + 1)).PasteSpecial(-4122)  # xlPasteFormats

Write-Output "PO Forecast sheet added; headers updated."
